$wb = $excel.ActiveWorkbook

# --- Sheet "Games": append new row 40 (Game 39), data for the MIA game that
#     has now been played (previously the first entry on the "Next" sheet). ---
$games = $wb.Worksheets.Item(1)

$games.Range("A40").Value = 39
$games.Range("B40").Value = 45306
$games.Range("B40").NumberFormat = $games.Range("B39").NumberFormat
$games.Range("C40").Value = -3
$games.Range("D40").Value = 95
$games.Range("E40").Value = 92.40000000000001
$games.Range("F40").Value = 0.404
$games.Range("G40").Value = 7
$games.Range("H40").Value = 15.8
$games.Range("I40").Value = 0.152
$games.Range("J40").Value = 93.09999999999999
$games.Range("K40").Value = "MIA"
$games.Range("L40").Value = 96
$games.Range("M40").Value = 0.425
$games.Range("N40").Value = 10.9
$games.Range("O40").Value = 12.8
$games.Range("P40").Value = 0.253
$games.Range("Q40").Value = 94.09999999999999
$games.Range("R40").Value = 1
$games.Range("S40").Value = 0

# --- Sheet "Next": remove the game that has now been played (old row 2,
#     MIA on 45306) - everything else shifts up by one row. ---
$next = $wb.Worksheets.Item(2)
$next.Rows.Item(2).Delete()
